# Update "想去人数" (interested-people count) values on the sheets that
# hold the conference data: "展览" (sheet1) and "全部类型" (sheet4).
# F3: 4187 -> 4189
# F4: 117  -> 118
# F5: 760  -> 761

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 4189
    $ws.Range("F4").Value = 118
    $ws.Range("F5").Value = 761
}
